# Update "想去人数" (interested count) values in column F for the sheets
# that contain this dataset: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1691
    4  = 780
    7  = 11847
    10 = 472
    11 = 402
    14 = 13444
    15 = 13354
    17 = 149
    23 = 156
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
